# Auto-generated edit script applying the commit diff
# 'Fixed results for 57 in journal'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("DCT")
$ws.Range("C2").Value = 0.4975
$ws.Range("D2").Value = 0.6639919759277833
$ws.Range("E2").Value = 0.4987443495730788
$ws.Range("F2").Value = 0.993
$ws.Range("G2").Value = 0.998
$ws.Range("C3").Value = 0.483
$ws.Range("D3").Value = 0.6492537313432836
$ws.Range("E3").Value = 0.4912731006160164
$ws.Range("F3").Value = 0.957
$ws.Range("G3").Value = 0.991
$ws.Range("C4").Value = 0.4975
$ws.Range("D4").Value = 0.6639919759277833
$ws.Range("E4").Value = 0.4987443495730788
$ws.Range("F4").Value = 0.993
$ws.Range("G4").Value = 0.998
$ws.Range("C5").Value = 0.4985
$ws.Range("D5").Value = 0.6639865996649916
$ws.Range("E5").Value = 0.4992443324937028
$ws.Range("F5").Value = 0.991
$ws.Range("G5").Value = 0.994
$ws.Range("C6").Value = 0.634
$ws.Range("D6").Value = 0.7306843267108168
$ws.Range("E6").Value = 0.5779976717112922
$ws.Range("F6").Value = 0.993
$ws.Range("G6").Value = 0.725
$ws.Range("C7").Value = 0.4995
$ws.Range("D7").Value = 0.6639812017455521
$ws.Range("E7").Value = 0.4997473471450227
$ws.Range("F7").Value = 0.989
$ws.Range("G7").Value = 0.99
$ws.Range("C8").Value = 0.4995
$ws.Range("D8").Value = 0.6657762938230384
$ws.Range("E8").Value = 0.4997493734335839
$ws.Range("F8").Value = 0.997
$ws.Range("G8").Value = 0.998
$ws.Range("C9").Value = 0.498
$ws.Range("D9").Value = 0.6639892904953146
$ws.Range("E9").Value = 0.4989939637826962
$ws.Range("F9").Value = 0.992
$ws.Range("G9").Value = 0.996
$ws.Range("C10").Value = 0.5065
$ws.Range("D10").Value = 0.6682352941176469
$ws.Range("E10").Value = 0.5032911392405063
$ws.Range("F10").Value = 0.994
$ws.Range("G10").Value = 0.981
$ws.Range("C11").Value = 0.9305
$ws.Range("D11").Value = 0.9253891572732152
$ws.Range("E11").Value = 0.9988412514484357
$ws.Range("F11").Value = 0.862
$ws.Range("G11").Value = 0.001

$ws = $wb.Worksheets.Item("GNB")
$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = 0.6666666666666666
$ws.Range("E3").Value = 0.5
$ws.Range("G3").Value = 1
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = 0.6666666666666666
$ws.Range("E4").Value = 0.5
$ws.Range("G4").Value = 1
$ws.Range("C5").Value = 0.9845
$ws.Range("D5").Value = 0.9847365829640571
$ws.Range("E5").Value = 0.9699321047526673
$ws.Range("G5").Value = 0.031
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 0.6666666666666666
$ws.Range("E6").Value = 0.5
$ws.Range("G6").Value = 1
$ws.Range("C7").Value = 0.9885
$ws.Range("D7").Value = 0.9886307464162135
$ws.Range("E7").Value = 0.9775171065493646
$ws.Range("G7").Value = 0.023
$ws.Range("C8").Value = 0.5
$ws.Range("D8").Value = 0.6666666666666666
$ws.Range("E8").Value = 0.5
$ws.Range("G8").Value = 1
$ws.Range("C9").Value = 0.958
$ws.Range("D9").Value = 0.9596928982725529
$ws.Range("E9").Value = 0.922509225092251
$ws.Range("G9").Value = 0.08400000000000001
$ws.Range("C10").Value = 0.9035
$ws.Range("D10").Value = 0.9119927040583675
$ws.Range("E10").Value = 0.8382229673093042
$ws.Range("G10").Value = 0.193
$ws.Range("C11").Value = 0.524
$ws.Range("D11").Value = 0.6775067750677507
$ws.Range("E11").Value = 0.5122950819672131
$ws.Range("G11").Value = 0.952

$ws = $wb.Worksheets.Item("KNN")
$ws.Range("C2").Value = 0.782
$ws.Range("D2").Value = 0.7212276214833758
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5639999999999999
$ws.Range("G2").Value = 0
$ws.Range("C3").Value = 0.8725000000000001
$ws.Range("D3").Value = 0.8843537414965986
$ws.Range("E3").Value = 0.8091286307053942
$ws.Range("F3").Value = 0.975
$ws.Range("G3").Value = 0.23
$ws.Range("C4").Value = 0.9495
$ws.Range("D4").Value = 0.94681411269089
$ws.Range("F4").Value = 0.899
$ws.Range("C5").Value = 0.9370000000000001
$ws.Range("D5").Value = 0.9327641408751334
$ws.Range("F5").Value = 0.874
$ws.Range("C6").Value = 0.9825
$ws.Range("D6").Value = 0.9821882951653944
$ws.Range("F6").Value = 0.965
$ws.Range("C7").Value = 0.982
$ws.Range("D7").Value = 0.9816700610997963
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.964
$ws.Range("G7").Value = 0
$ws.Range("C8").Value = 0.954
$ws.Range("D8").Value = 0.9517819706498952
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.908
$ws.Range("G8").Value = 0
$ws.Range("C9").Value = 0.9325
$ws.Range("D9").Value = 0.9276139410187667
$ws.Range("F9").Value = 0.865
$ws.Range("C10").Value = 0.9835
$ws.Range("D10").Value = 0.9832231825114387
$ws.Range("F10").Value = 0.967
$ws.Range("C11").Value = 0.984
$ws.Range("D11").Value = 0.983739837398374
$ws.Range("F11").Value = 0.968

$ws = $wb.Worksheets.Item("SVM_L")
$ws.Range("C2").Value = 0.5315
$ws.Range("D2").Value = 0.1185324553151458
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.063
$ws.Range("G2").Value = 0
$ws.Range("C3").Value = 0.5004999999999999
$ws.Range("D3").Value = 0.6668889629876625
$ws.Range("E3").Value = 0.5002501250625313
$ws.Range("G3").Value = 0.999
$ws.Range("C4").Value = 0.997
$ws.Range("D4").Value = 0.997002997002997
$ws.Range("E4").Value = 0.9960079840319361
$ws.Range("F4").Value = 0.998
$ws.Range("G4").Value = 0.004
$ws.Range("C5").Value = 0.9565
$ws.Range("D5").Value = 0.9545216936748563
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.913
$ws.Range("G5").Value = 0
$ws.Range("C6").Value = 0.5004999999999999
$ws.Range("D6").Value = 0.6666666666666667
$ws.Range("E6").Value = 0.5002503755633451
$ws.Range("G6").Value = 0.998
$ws.Range("C7").Value = 0.9985000000000001
$ws.Range("D7").Value = 0.99849774661993
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.997
$ws.Range("G7").Value = 0
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("C9").Value = 0.9955000000000001
$ws.Range("D9").Value = 0.9954796584630838
$ws.Range("F9").Value = 0.991
$ws.Range("C10").Value = 0.9935
$ws.Range("D10").Value = 0.9934574735782586
$ws.Range("F10").Value = 0.987
$ws.Range("C11").Value = 0.996
$ws.Range("D11").Value = 0.9959839357429718
$ws.Range("F11").Value = 0.992

$ws = $wb.Worksheets.Item("SVM-R")
$ws.Range("C2").Value = 0.548
$ws.Range("D2").Value = 0.1751824817518248
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.096
$ws.Range("G2").Value = 0
$ws.Range("C3").Value = 0.983
$ws.Range("D3").Value = 0.983184965380811
$ws.Range("E3").Value = 0.9726027397260274
$ws.Range("F3").Value = 0.994
$ws.Range("G3").Value = 0.028
$ws.Range("C4").Value = 0.972
$ws.Range("D4").Value = 0.971252566735113
$ws.Range("E4").Value = 0.9978902953586498
$ws.Range("F4").Value = 0.946
$ws.Range("G4").Value = 0.002
$ws.Range("C5").Value = 0.948
$ws.Range("D5").Value = 0.9451476793248946
$ws.Range("F5").Value = 0.896
$ws.Range("C6").Value = 0.995
$ws.Range("D6").Value = 0.9950199203187251
$ws.Range("E6").Value = 0.9910714285714286
$ws.Range("F6").Value = 0.999
$ws.Range("G6").Value = 0.008999999999999999
$ws.Range("C7").Value = 0.9965000000000001
$ws.Range("D7").Value = 0.9964877069744105
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.993
$ws.Range("G7").Value = 0
$ws.Range("C8").Value = 0.9995000000000001
$ws.Range("D8").Value = 0.9994997498749374
$ws.Range("F8").Value = 0.999
$ws.Range("C9").Value = 0.993
$ws.Range("D9").Value = 0.9929506545820745
$ws.Range("F9").Value = 0.986
$ws.Range("C10").Value = 0.993
$ws.Range("D10").Value = 0.9929506545820745
$ws.Range("F10").Value = 0.986
$ws.Range("C11").Value = 0.993
$ws.Range("D11").Value = 0.9929506545820745
$ws.Range("F11").Value = 0.986

$ws = $wb.Worksheets.Item("Training Time")
$ws.Range("C2").Value = 66.77434110641479
$ws.Range("C3").Value = 3.152542352676392
$ws.Range("C4").Value = 0.8976325988769531
$ws.Range("C5").Value = 143.1019971370697
$ws.Range("C6").Value = 2537.503065347672
